$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "Create a Game Character Face from a Single Portrait!`nhttps://www.catalyzex.com/p"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2706"

$ws.Range("D25").Value = "[바람돌이/딥러닝] GAN(3) - DCGAN(Deep Convolutional Generative Adversarial Networks) 논문 이론 및 리뷰"
$ws.Range("E25").Value = "https://blog.naver.com/winddori2002/222239283526"
